# "added 4wk low sales check"
# Update PO_Forecast values in column B to reflect the new 4-week low sales check.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = 47
    3  = 44
    4  = 42
    5  = 35
    7  = 29
    10 = 22
    11 = 20
    12 = 17
    13 = 15
    14 = 13
    19 = 31
    22 = 44
    27 = 160
    29 = 181
    30 = 192
    32 = 213
    33 = 224
    34 = 234
    35 = 245
    36 = 256
    37 = 266
    38 = 277
    39 = 287
}

foreach ($row in $updates.Keys) {
    $ws.Range("B$row").Value = $updates[$row]
}
